$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet4")

# -----------------------------------------------------------------
# The IP_Center_Specific_Code column (column C) on Sheet1 held raw
# "center_nameN" values. Mask/rename them to the generic
# "hms batch name N" scheme (the data column only - row 10 is blank
# and is skipped, exactly like the source edit).
# -----------------------------------------------------------------
$ws1.Range("C3").Value  = "hms batch name 1"
$ws1.Range("C4").Value  = "hms batch name 1a"
$ws1.Range("C5").Value  = "hms batch name 2a"
$ws1.Range("C6").Value  = "hms batch name 2"
$ws1.Range("C7").Value  = "hms batch name 3c"
$ws1.Range("C8").Value  = "hms batch name 3-3"
$ws1.Range("C9").Value  = "hms batch name 7"
$ws1.Range("C11").Value = "hms batch name 9"
$ws1.Range("C12").Value = "hms batch name 10"
$ws1.Range("C13").Value = "hms batch name 11"
$ws1.Range("C14").Value = "hms batch name 12"
$ws1.Range("C15").Value = "hms batch name 13"

# -----------------------------------------------------------------
# Reflect the editor's selection: the edited column-C cells
# (C3:C9 and C11:C15, skipping the blank C10) are left selected on
# Sheet1, and Sheet4 keeps A1 active while carrying the same
# multi-range selection context.
# -----------------------------------------------------------------
$ws1.Activate()
$editedRange = $ws1.Range("C3:C9,C11:C15")
$editedRange.Select() | Out-Null

$ws2.Activate()
$ws2.Range("A1").Select() | Out-Null
$ws1.Activate() | Out-Null
